$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Id=1 record - bump "Num of changes" and "Update Time"
$ws.Cells.Item(2, 3).Value = "'19"
$ws.Cells.Item(2, 3).ClearFormats()
$ws.Cells.Item(2, 4).Value = '"2023-06-10T19:42:36.575Z"'

# Row 3: Id=2 record - fill in "Is Manual", reset "Num of changes", update "Update Time"
$ws.Cells.Item(3, 2).Value = "'true"
$ws.Cells.Item(3, 2).ClearFormats()
$ws.Cells.Item(3, 3).Value = "'3"
$ws.Cells.Item(3, 3).ClearFormats()
$ws.Cells.Item(3, 4).Value = '"2023-06-10T19:43:34.378Z"'

# Row 4: Id=3 record - reset "Num of changes", update "Update Time"
$ws.Cells.Item(4, 3).Value = "'3"
$ws.Cells.Item(4, 3).ClearFormats()
$ws.Cells.Item(4, 4).Value = '"2023-06-10T19:53:08.505Z"'

# Row 5: brand-new Id=4 record
$ws.Cells.Item(5, 1).Value = "'4"
$ws.Cells.Item(5, 1).ClearFormats()
$ws.Cells.Item(5, 2).Value = "'true"
$ws.Cells.Item(5, 2).ClearFormats()
$ws.Cells.Item(5, 3).Value = "'6"
$ws.Cells.Item(5, 3).ClearFormats()
$ws.Cells.Item(5, 4).Value = '"2023-06-10T19:57:51.980Z"'
